# Restore the "USP18" data column that was accidentally dropped from the
# ImmGen signaling workbook: append header "USP18" in Z1 and the matching
# per-cell-type values in Z2:Z20 (rows already present in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 26).Value = "USP18"

$ws.Cells.Item(2, 26).Value = 308.851
$ws.Cells.Item(3, 26).Value = 299.661
$ws.Cells.Item(4, 26).Value = 115.9
$ws.Cells.Item(5, 26).Value = 90.2267
$ws.Cells.Item(6, 26).Value = 120.015
$ws.Cells.Item(7, 26).Value = 207.27
$ws.Cells.Item(8, 26).Value = 179.981
$ws.Cells.Item(9, 26).Value = 178.26
$ws.Cells.Item(10, 26).Value = 78.4159
$ws.Cells.Item(11, 26).Value = 132.335
$ws.Cells.Item(12, 26).Value = 167.415
$ws.Cells.Item(13, 26).Value = 74.055
$ws.Cells.Item(14, 26).Value = 146.105
$ws.Cells.Item(15, 26).Value = 152.443
$ws.Cells.Item(16, 26).Value = 147.382
$ws.Cells.Item(17, 26).Value = 223.252
$ws.Cells.Item(18, 26).Value = 255.253
$ws.Cells.Item(19, 26).Value = 126.112
$ws.Cells.Item(20, 26).Value = 136.67

# Match the reviewed workbook's view state: scrolled so column K is
# leftmost, with the newly restored Z2 cell selected.
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("Z2").Select()
